# Auto-generated script to refresh cached market-board values
# as captured by the scheduled runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 1738.5
$ws.Cells.Item(96, 9).Value = 983.625
$ws.Cells.Item(96, 10).Value = 2342.4
$ws.Cells.Item(96, 11).Value = 2950.875
$ws.Cells.Item(96, 12).Value = 7027.200000000001
$ws.Cells.Item(96, 13).Value = -1577.875
$ws.Cells.Item(96, 14).Value = -9773.200000000001
$ws.Cells.Item(116, 8).Value = 2157.25
$ws.Cells.Item(116, 9).Value = 1782
$ws.Cells.Item(116, 10).Value = 2327.818
$ws.Cells.Item(116, 11).Value = 1782
$ws.Cells.Item(116, 12).Value = 2327.818
$ws.Cells.Item(116, 13).Value = 1660
$ws.Cells.Item(116, 14).Value = -9211.817999999999
$ws.Cells.Item(127, 8).Value = 3983.8
$ws.Cells.Item(127, 9).Value = 871.75
$ws.Cells.Item(127, 10).Value = 6058.5
$ws.Cells.Item(127, 11).Value = 2615.25
$ws.Cells.Item(127, 12).Value = 18175.5
$ws.Cells.Item(127, 13).Value = 2344.75
$ws.Cells.Item(127, 14).Value = -28095.5
$ws.Cells.Item(132, 8).Value = 2166427.8
$ws.Cells.Item(132, 9).Value = 2343769
$ws.Cells.Item(132, 10).Value = 2865.2
$ws.Cells.Item(132, 11).Value = 7031307
$ws.Cells.Item(132, 12).Value = 8595.599999999999
$ws.Cells.Item(132, 13).Value = -7028777
$ws.Cells.Item(132, 14).Value = -13655.6
$ws.Cells.Item(137, 8).Value = 5406553
$ws.Cells.Item(137, 9).Value = 1160.1666
$ws.Cells.Item(137, 10).Value = 10527452
$ws.Cells.Item(137, 11).Value = 3480.4998
$ws.Cells.Item(137, 12).Value = 31582356
$ws.Cells.Item(137, 13).Value = -930.4998000000001
$ws.Cells.Item(137, 14).Value = -31587456
$ws.Cells.Item(141, 8).Value = 2393.6365
$ws.Cells.Item(141, 9).Value = 1416.25
$ws.Cells.Item(141, 10).Value = 5000
$ws.Cells.Item(141, 11).Value = 4248.75
$ws.Cells.Item(141, 12).Value = 15000
$ws.Cells.Item(141, 13).Value = 931.25
$ws.Cells.Item(141, 14).Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 16485.494
$ws.Cells.Item(32, 9).Value = 16969.184
$ws.Cells.Item(32, 10).Value = 7900
$ws.Cells.Item(32, 11).Value = 16969.184
$ws.Cells.Item(32, 12).Value = 7900
$ws.Cells.Item(32, 13).Value = -16682.184
$ws.Cells.Item(32, 14).Value = -8474
$ws.Cells.Item(61, 8).Value = 929.2456
$ws.Cells.Item(61, 9).Value = 650.1163
$ws.Cells.Item(61, 10).Value = 1786.5714
$ws.Cells.Item(61, 11).Value = 650.1163
$ws.Cells.Item(61, 12).Value = 1786.5714
$ws.Cells.Item(61, 13).Value = -438.1163
$ws.Cells.Item(61, 14).Value = -2210.5714
$ws.Cells.Item(63, 8).Value = 669086.3
$ws.Cells.Item(63, 9).Value = 911127.75
$ws.Cells.Item(63, 10).Value = 3472.5
$ws.Cells.Item(63, 11).Value = 911127.75
$ws.Cells.Item(63, 12).Value = 3472.5
$ws.Cells.Item(63, 13).Value = -910441.75
$ws.Cells.Item(63, 14).Value = -4844.5
$ws.Cells.Item(66, 8).Value = 669086.3
$ws.Cells.Item(66, 9).Value = 911127.75
$ws.Cells.Item(66, 10).Value = 3472.5
$ws.Cells.Item(66, 11).Value = 4555638.75
$ws.Cells.Item(66, 12).Value = 17362.5
$ws.Cells.Item(66, 13).Value = -4552206.75
$ws.Cells.Item(66, 14).Value = -24226.5
$ws.Cells.Item(74, 8).Value = 588.53845
$ws.Cells.Item(74, 9).Value = 445.54544
$ws.Cells.Item(74, 11).Value = 445.54544
$ws.Cells.Item(74, 13).Value = 428.45456
$ws.Cells.Item(77, 8).Value = 588.53845
$ws.Cells.Item(77, 9).Value = 445.54544
$ws.Cells.Item(77, 11).Value = 2227.7272
$ws.Cells.Item(77, 13).Value = 2140.2728
$ws.Cells.Item(80, 8).Value = 24233.334
$ws.Cells.Item(80, 10).Value = 24233.334
$ws.Cells.Item(80, 12).Value = 24233.334
$ws.Cells.Item(80, 14).Value = -26229.334
$ws.Cells.Item(83, 8).Value = 24233.334
$ws.Cells.Item(83, 10).Value = 24233.334
$ws.Cells.Item(83, 12).Value = 72700.00199999999
$ws.Cells.Item(83, 14).Value = -82684.00199999999
$ws.Cells.Item(102, 8).Value = 1603.6666
$ws.Cells.Item(102, 9).Value = 1603.6666
$ws.Cells.Item(102, 11).Value = 1603.6666
$ws.Cells.Item(102, 13).Value = 18.33339999999998
$ws.Cells.Item(122, 8).Value = 13180.223
$ws.Cells.Item(122, 9).Value = 2230.2856
$ws.Cells.Item(122, 11).Value = 6690.8568
$ws.Cells.Item(122, 13).Value = -4240.8568
$ws.Cells.Item(136, 8).Value = 929.2456
$ws.Cells.Item(136, 9).Value = 650.1163
$ws.Cells.Item(136, 10).Value = 1786.5714
$ws.Cells.Item(136, 11).Value = 1950.3489
$ws.Cells.Item(136, 12).Value = 5359.7142
$ws.Cells.Item(136, 13).Value = 599.6511
$ws.Cells.Item(136, 14).Value = -10459.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 366.44446
$ws.Cells.Item(80, 9).Value = 822
$ws.Cells.Item(80, 10).Value = 191.23077
$ws.Cells.Item(80, 11).Value = 822
$ws.Cells.Item(80, 12).Value = 191.23077
$ws.Cells.Item(80, 13).Value = 176
$ws.Cells.Item(80, 14).Value = -2187.23077
$ws.Cells.Item(82, 8).Value = 42172.383
$ws.Cells.Item(82, 10).Value = 65651.625
$ws.Cells.Item(82, 12).Value = 65651.625
$ws.Cells.Item(82, 14).Value = -66417.625
$ws.Cells.Item(83, 8).Value = 366.44446
$ws.Cells.Item(83, 9).Value = 822
$ws.Cells.Item(83, 10).Value = 191.23077
$ws.Cells.Item(83, 11).Value = 4110
$ws.Cells.Item(83, 12).Value = 956.15385
$ws.Cells.Item(83, 13).Value = 882
$ws.Cells.Item(83, 14).Value = -10940.15385
$ws.Cells.Item(85, 8).Value = 42172.383
$ws.Cells.Item(85, 10).Value = 65651.625
$ws.Cells.Item(85, 12).Value = 65651.625
$ws.Cells.Item(85, 14).Value = -68303.625
$ws.Cells.Item(134, 8).Value = 20286.314
$ws.Cells.Item(134, 9).Value = 25415.666
$ws.Cells.Item(134, 10).Value = 2333.5833
$ws.Cells.Item(134, 11).Value = 76246.99800000001
$ws.Cells.Item(134, 12).Value = 7000.749899999999
$ws.Cells.Item(134, 13).Value = -73711.99800000001
$ws.Cells.Item(134, 14).Value = -12070.7499

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 951.4286
$ws.Cells.Item(18, 9).Value = 300
$ws.Cells.Item(18, 10).Value = 2124
$ws.Cells.Item(18, 11).Value = 900
$ws.Cells.Item(18, 12).Value = 6372
$ws.Cells.Item(18, 13).Value = -731
$ws.Cells.Item(18, 14).Value = -6710
$ws.Cells.Item(68, 8).Value = 1361.2858
$ws.Cells.Item(68, 9).Value = 1292.7894
$ws.Cells.Item(68, 10).Value = 1417.8695
$ws.Cells.Item(68, 11).Value = 3878.3682
$ws.Cells.Item(68, 12).Value = 4253.6085
$ws.Cells.Item(68, 13).Value = -3067.3682
$ws.Cells.Item(68, 14).Value = -5875.6085
$ws.Cells.Item(71, 8).Value = 1361.2858
$ws.Cells.Item(71, 9).Value = 1292.7894
$ws.Cells.Item(71, 10).Value = 1417.8695
$ws.Cells.Item(71, 11).Value = 11635.1046
$ws.Cells.Item(71, 12).Value = 12760.8255
$ws.Cells.Item(71, 13).Value = -7579.104599999999
$ws.Cells.Item(71, 14).Value = -20872.8255
$ws.Cells.Item(94, 8).Value = 8106.75
$ws.Cells.Item(94, 9).Value = 2900
$ws.Cells.Item(94, 10).Value = 9842.333000000001
$ws.Cells.Item(94, 11).Value = 8700
$ws.Cells.Item(94, 12).Value = 29526.999
$ws.Cells.Item(94, 13).Value = -8024
$ws.Cells.Item(94, 14).Value = -30878.999
$ws.Cells.Item(131, 8).Value = 3177786.2
$ws.Cells.Item(131, 10).Value = 5821123
$ws.Cells.Item(131, 12).Value = 17463369
$ws.Cells.Item(131, 14).Value = -17473449

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 11338231
$ws.Cells.Item(70, 9).Value = 17004454
$ws.Cells.Item(70, 10).Value = 5785.3335
$ws.Cells.Item(70, 11).Value = 17004454
$ws.Cells.Item(70, 12).Value = 5785.3335
$ws.Cells.Item(70, 13).Value = -17004184
$ws.Cells.Item(70, 14).Value = -6325.3335
$ws.Cells.Item(73, 8).Value = 11338231
$ws.Cells.Item(73, 9).Value = 17004454
$ws.Cells.Item(73, 10).Value = 5785.3335
$ws.Cells.Item(73, 11).Value = 17004454
$ws.Cells.Item(73, 12).Value = 5785.3335
$ws.Cells.Item(73, 13).Value = -17003518
$ws.Cells.Item(73, 14).Value = -7657.3335
$ws.Cells.Item(80, 8).Value = 4238.25
$ws.Cells.Item(80, 9).Value = 3080.4167
$ws.Cells.Item(80, 10).Value = 5975
$ws.Cells.Item(80, 11).Value = 3080.4167
$ws.Cells.Item(80, 12).Value = 5975
$ws.Cells.Item(80, 13).Value = -2082.4167
$ws.Cells.Item(80, 14).Value = -7971
$ws.Cells.Item(83, 8).Value = 4238.25
$ws.Cells.Item(83, 9).Value = 3080.4167
$ws.Cells.Item(83, 10).Value = 5975
$ws.Cells.Item(83, 11).Value = 15402.0835
$ws.Cells.Item(83, 12).Value = 29875
$ws.Cells.Item(83, 13).Value = -10410.0835
$ws.Cells.Item(83, 14).Value = -39859
$ws.Cells.Item(99, 8).Value = 3752.1667
$ws.Cells.Item(99, 9).Value = 1502.6
$ws.Cells.Item(99, 11).Value = 1502.6
$ws.Cells.Item(99, 13).Value = 743.4000000000001
$ws.Cells.Item(113, 8).Value = 16668444
$ws.Cells.Item(113, 9).Value = 35715244
$ws.Cells.Item(113, 10).Value = 2493.75
$ws.Cells.Item(113, 11).Value = 35715244
$ws.Cells.Item(113, 12).Value = 2493.75
$ws.Cells.Item(113, 13).Value = -35713074
$ws.Cells.Item(113, 14).Value = -6833.75
$ws.Cells.Item(122, 8).Value = 7276.3335
$ws.Cells.Item(122, 9).Value = 20100
$ws.Cells.Item(122, 11).Value = 60300
$ws.Cells.Item(122, 13).Value = -57850
$ws.Cells.Item(132, 8).Value = 77006.81
$ws.Cells.Item(132, 9).Value = 113603.89
$ws.Cells.Item(132, 10).Value = 3812.6667
$ws.Cells.Item(132, 11).Value = 340811.67
$ws.Cells.Item(132, 12).Value = 11438.0001
$ws.Cells.Item(132, 13).Value = -338281.67
$ws.Cells.Item(132, 14).Value = -16498.0001
$ws.Cells.Item(135, 8).Value = 48875
$ws.Cells.Item(135, 10).Value = 48875
$ws.Cells.Item(135, 12).Value = 48875
$ws.Cells.Item(135, 14).Value = -59015

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 760
$ws.Cells.Item(22, 9).Value = 700
$ws.Cells.Item(22, 10).Value = 800
$ws.Cells.Item(22, 11).Value = 700
$ws.Cells.Item(22, 12).Value = 800
$ws.Cells.Item(22, 13).Value = -405
$ws.Cells.Item(22, 14).Value = -1390
$ws.Cells.Item(27, 8).Value = 760
$ws.Cells.Item(27, 9).Value = 700
$ws.Cells.Item(27, 10).Value = 800
$ws.Cells.Item(27, 11).Value = 700
$ws.Cells.Item(27, 12).Value = 800
$ws.Cells.Item(27, 13).Value = -593
$ws.Cells.Item(27, 14).Value = -1014
$ws.Cells.Item(132, 8).Value = 16247.4
$ws.Cells.Item(132, 9).Value = 24767.223
$ws.Cells.Item(132, 10).Value = 3467.6667
$ws.Cells.Item(132, 11).Value = 74301.66900000001
$ws.Cells.Item(132, 12).Value = 10403.0001
$ws.Cells.Item(132, 13).Value = -71771.66900000001
$ws.Cells.Item(132, 14).Value = -15463.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1156.0294
$ws.Cells.Item(132, 9).Value = 978.3125
$ws.Cells.Item(132, 11).Value = 2934.9375
$ws.Cells.Item(132, 13).Value = -404.9375

